# Insert a new weekly data row at row 169 (Feria Lagunitas de Puerto Montt -
# Betarraga), shifting the existing rows 169-207 down to 170-208 and
# extending the used range from A1:R207 to A1:R208.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("169:169").Insert()

$ws.Range("A169").Value = 4
$ws.Range("B169").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C169").Value = "Los Lagos"
$ws.Range("D169").Value = 44511
$ws.Range("E169").Value = 10
$ws.Range("F169").Value = 100114014
$ws.Range("G169").Value = "Betarraga"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 500
$ws.Range("K169").Value = 1000
$ws.Range("L169").Value = 1000
$ws.Range("M169").Value = 1000
$ws.Range("N169").Value = "$/paquete 5 unidades"
$ws.Range("O169").Value = "Región del Maule"
$ws.Range("P169").Value = 200
$ws.Range("Q169").Value = 5
$ws.Range("R169").Value = "Hortaliza"
